$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 212
$ws.Range("I2").Value = 543
$ws.Range("J2").Value = 2269
$ws.Range("K2").Value = 11
$ws.Range("L2").Value = 612
$ws.Range("M2").Value = 40
$ws.Range("N2").Value = 379
$ws.Range("P2").Value = 7
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 18
$ws.Range("S2").Value = 239
$ws.Range("T2").Value = 400
$ws.Range("U2").Value = 27
$ws.Range("V2").Value = 3746
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 3657
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 56
$ws.Range("AA2").Value = 22
